# "fix: add two drinks"
# - Sheet2: correct two mis-typed drink names ("唐师傅..." -> "康师傅...")
# - Sheet1: append two new drink rows at the bottom of the list

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: fix the two drink names -------------------------------------
$ws2.Activate()
$ws2.Range("B41").Value = "康师傅水蜜桃"
$ws2.Range("B38").Value = "康师傅绿茶"
$ws2.Range("B1:B2").Select()

# --- Sheet1: add the two new drinks at rows 99-100 -----------------------
$ws1.Activate()

$ws1.Range("B99").Value = "冰公主猕猴桃汁"
$ws1.Range("C99").Value = "一个商品介绍"
$ws1.Range("D99").Value = 6

$ws1.Range("B100").Value = "百加可植物饮料"
$ws1.Range("C100").Value = "一个商品介绍"
$ws1.Range("D100").Value = 7

$ws1.Range("D102").Select()

$wb.Save()
